$p = $ppt.ActivePresentation

# --- Slide 2: "Conway's Game of Life" -> "Conway's Algorithm" ---
$s2 = $p.Slides.Item(2)
$title2 = $s2.Shapes.Item(1).TextFrame.TextRange
# Replace the "Game of Life" portion (chars 10-21) with "Algorithm",
# leaving "Conway’s " (chars 1-9) in its own run untouched.
$title2.Characters(10, 12).Text = "Algorithm"

# --- Slide 3: "Conway's Game of Life" -> "Conway's Algorithm - Example" ---
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1).TextFrame.TextRange
$title3.Characters(10, 12).Text = "Algorithm - Example"
